# Fix lỗi mã số thuế - bị nối thêm Tên công ty vào
# Populate the "TINH BINH DINH" sheet with the corrected tax-code lookup rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TINH BINH DINH")

# Row 1
$ws.Range("A1").Value = 4101598873
$ws.Range("B1").Value = 44414
$ws.Range("C1").Value = "CÔNG TY TNHH VẬN TẢI TUẤN MINH BÌNH ĐỊNH"
$ws.Range("D1").Value = "Lô 15-16 Khu Đô thị An Phú Thịnh, Phường Đống Đa, Thành phố Quy Nhơn, Tỉnh Bình Định"
$ws.Range("E1").Value = 965972999
$ws.Range("F1").Value = "Nguyễn Quy Khoa"

# Row 2
$ws.Range("A2").Value = 4101598880
$ws.Range("B2").Value = 44414
$ws.Range("C2").Value = "CÔNG TY TNHH THƯƠNG MẠI - TỔNG HỢP BÌNH VƯƠNG"
$ws.Range("D2").Value = "Số 295 Nguyễn Thị Minh Khai, Phường Nguyễn Văn Cừ, Thành phố Quy Nhơn, Tỉnh Bình Định"
$ws.Range("E2").Value = 963555405
$ws.Range("F2").Value = "Bùi Quốc Thắng"

# Re-use the workbook's existing short-date cell style (already present on
# "Sheet1"!B24:B50) for the date column instead of letting Excel mint a new,
# duplicate number format.
$wsDates = $wb.Worksheets.Item("Sheet1")
$wsDates.Range("B24").Copy()
$ws.Range("B1:B2").PasteSpecial(-4122)
